$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AppControl")

# Fill the "NA" values for B30, B31, B32, B33 (fixes null pointer exception
# when these cells were previously blank).
$ws.Range("B30").Value = "NA"
$ws.Range("B31").Value = "NA"
$ws.Range("B32").Value = "NA"
$ws.Range("B33").Value = "NA"

# Update the sheet view: scroll so row 28 is at top and select B36.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("B36").Select() | Out-Null
